# Adding ACTIVATE cases to library:
# Insert a new row above the "area_fraction_cover_of_liquid_cloud" entry (row 21)
# for the new "cloud_area_fraction" / "clt" variable, copying the row format
# from the row above it (row 20) and filling in the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21 (pushes existing row 21+ down by one)
$ws.Range("A21:Z21").EntireRow.Insert()

# Copy the formatting from row 20 (the row just above the new one) into the
# newly inserted row 21, so the new row matches the table's look & feel.
$ws.Range("A20:Z20").Copy()
$ws.Range("A21:Z21").PasteSpecial(-4122)
$ws.Rows.Item(21).RowHeight = 15.75

# Populate the new row with the new "cloud_area_fraction" entry.
$ws.Range("A21").Value = "cloud_area_fraction"
$ws.Range("B21").Value = "clt"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "time"
$ws.Range("E21").Value = "diagnosed cloud cover"

# Update the active selection to reflect where the editor left off.
$ws.Range("E21").Select()
